# Remove the leftover full-slide "TBD" placeholder rectangle (id=5,
# name "직사각형 4") from slide 280 (the second slide), making the
# slide's real content the default/visible screen again.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Id -eq 5) {
        $sh.Delete()
    }
}
